$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "H 72" record) - all data below shifts up by one row
$ws.Rows.Item(2).Delete()
